$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '29.209.81'
$ws.Range("E2").Value = '  -0.93%  '

# Row 3
$ws.Range("D3").Value = '1.859.09'
$ws.Range("E3").Value = '  -0.63%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.12%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7152'
$ws.Range("E5").Value = '  -0.27%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '240.45'
$ws.Range("E6").Value = '  +0.59%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.000'
$ws.Range("E7").Value = '  -0.11%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07752'
$ws.Range("E8").Value = '  -0.50%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3076'
$ws.Range("E9").Value = '  +0.22%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '25.15'
$ws.Range("E10").Value = '  -0.33%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08252'
$ws.Range("E11").Value = '  +0.22%  '

# Row 12
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.852.40'
$ws.Range("E12").Value = '  -0.96%  '

# Row 13
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.241'
$ws.Range("E13").Value = '  +0.28%  '

# Row 14
$ws.Range("B14").Value = 'Polygon'
$ws.Range("C14").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7175'
$ws.Range("E14").Value = '  -0.37%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '90.23'
$ws.Range("E15").Value = '  +0.36%  '

# Row 16
$ws.Range("D16").Value = '29.181.75'
$ws.Range("E16").Value = '  -1.23%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.867'
$ws.Range("E17").Value = '  +0.78%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '244.23'
$ws.Range("E18").Value = '  +1.41%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007790'
$ws.Range("E19").Value = '  -0.74%  '

# Row 20
$ws.Range("E20").Value = '  -1.08%  '

# Row 21
$ws.Range("D21").Value = '2.106.44'
$ws.Range("E21").Value = '  -1.70%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9990'
$ws.Range("E22").Value = '  -0.17%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.982'
$ws.Range("E23").Value = '  +3.39%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.000'
$ws.Range("E24").Value = '  -0.14%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1597'
$ws.Range("E25").Value = '  +2.24%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '162.43'
$ws.Range("E26").Value = '  -0.11%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.922'
$ws.Range("E27").Value = '  -0.39%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.27'
$ws.Range("E28").Value = '  -0.22%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.496'
$ws.Range("E29").Value = '  +0.91%  '

# Row 30
$ws.Range("E30").Value = '  -3.16%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.401'
$ws.Range("E31").Value = '  +1.72%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.202'
$ws.Range("E32").Value = '  +3.16%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05190'
$ws.Range("E33").Value = '  -1.16%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.910'
$ws.Range("E34").Value = '  -1.14%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.171'
$ws.Range("E35").Value = '  -2.11%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7273'
$ws.Range("E36").Value = '  +1.61%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.677'
$ws.Range("E37").Value = '  +0.16%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01856'
$ws.Range("E38").Value = '  -0.53%  '

# Row 39
$ws.Range("E39").Value = '  -1.34%  '

# Row 40
$ws.Range("D40").Value = '1.149.48'
$ws.Range("E40").Value = '  -2.12%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9040'
$ws.Range("E41").Value = '  -0.21%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.145'
$ws.Range("E42").Value = '  +2.56%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '72.32'
$ws.Range("E43").Value = '  +1.40%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9990'
$ws.Range("E44").Value = '  -0.24%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '101.61'
$ws.Range("E45").Value = '  -0.67%  '

# Row 46
$ws.Range("D46").Value = '2.000.39'
$ws.Range("E46").Value = '  -1.82%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5230'
$ws.Range("E47").Value = '  -2.39%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.767'
$ws.Range("E48").Value = '  +0.34%  '

# Row 49
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00000000120'
$ws.Range("E49").Value = '  -1.37%  '

# Row 50
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.334'
$ws.Range("E50").Value = '  +2.07%  '

# Row 51
$ws.Range("E51").Value = '  +1.48%  '
